$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G header (row 1) - inline string "model4"
$ws.Range("G1").Value = "model4"

# Add G5 = 1 (numeric)
$ws.Range("G5").Value = 1

# Add G6 = 1 (numeric), F6 stays as 1
$ws.Range("G6").Value = 1
